$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating the "总计" sheet
#    (so it inherits the same header/row styling), then place it
#    right before "总计".
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($total)
$newQ = $wb.Worksheets.Item("总计 (2)")
$newQ.Name = "2022-Q1"

# Re-fetch "总计" since its Index/handle can go stale after the sheet move.
$total = $wb.Worksheets.Item("总计")

# Drop the old data rows (3..6), keep row 1 (header) and row 2 (which
# will become the single fund data row) so their formatting carries over.
$newQ.Rows("3:6").Delete()
$newQ.Range("B2:D2").ClearContents()

# Extend the header formatting (style of D1) across the new columns E1:H1.
$newQ.Range("D1").Copy()
$newQ.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$newQ.Range("A1").Copy() | Out-Null

# Header row text
$newQ.Range("B1").Value = "基金代码"
$newQ.Range("C1").Value = "基金名称"
$newQ.Range("D1").Value = "基金规模"
$newQ.Range("E1").Value = "股票总仓位"
$newQ.Range("F1").Value = "仓位占比"
$newQ.Range("G1").Value = "持有市值(亿元)"
$newQ.Range("H1").Value = "仓位排名"

# Single fund data row
$newQ.Range("A2").Value = 0
$newQ.Range("B2").Value = "'206009"
$newQ.Range("C2").Value = "鹏华新兴产业混合"
$newQ.Range("D2").Value = "'44.95"
$newQ.Range("E2").Value = "'90.17"
$newQ.Range("F2").Value = "'5.15"
$newQ.Range("G2").Value = "'2.3149"
$newQ.Range("H2").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top data row for 2022-Q1
#    and renumber the existing index column.
# ------------------------------------------------------------------
$total.Rows("2:2").Insert()

# Re-apply the index-column style (copied from the row pushed down to 3)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122) | Out-Null
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "'2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 2.31
